# Applies the "new version with timestamp" update:
#  1. Removes the item row for "املاح افونا" (row 18), shifting all rows below it up by one.
#  2. Updates the recalculated total (previously 985.44 -> 960.44).
#  3. Refreshes the generation timestamp in the footer (11:15 AM -> 11:21 AM).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the item row for "املاح افونا" (row 18) - shifts rows 19-26 up to 18-25.
$ws.Rows.Item(18).Delete()

# 2. Re-number the serial-number column ("م") for the remaining item rows (18-23),
#    since the source report regenerates this sequential index (1..17) rather than
#    keeping the old numbers shifted with the rows.
for ($r = 18; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 6
}

# 3. Fix the recalculated grand-total cell (now on row 24, column P), restoring the
#    taller row height that belonged to the deleted item row (the shifted-up total
#    row otherwise keeps its own, shorter height).
$ws.Cells.Item(24, 16).Value = 960.44
$ws.Rows.Item(24).RowHeight = 25.5

# 4. Refresh the footer timestamp (now on row 25, column A).
$ws.Cells.Item(25, 1).Value = "Tuesday, 30 September, 2025 11:21 AM"
